# Updated cryptos list - apply latest price/volume(1h) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.528.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.444.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.445.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +2.63%  "
$ws.Range("E11").Value = "  +1.20%  "
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.032.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.88%  "
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.444.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.609.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.78%  "
$ws.Range("E19").Value = "  +8.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "
$ws.Range("E21").Value = "  +0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.587.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.53"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.27"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.37%  "
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("E36").Value = "  +1.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.90%  "
$ws.Range("E39").Value = "  +0.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.59"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.612.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.59%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.31%  "
